# Updates the cryptocurrency price/volume table on Sheet1 with freshly
# scraped values (coinranking.com), as produced by the scheduled
# "Updated cryptos list ... with GitHub Actions" job.
#
# Every cell in this table - including columns D (Price) and E
# (Volume(1h)), and for the few rows whose ranking shuffled, columns B
# (Coin) and C (Link) - is stored as plain text in the workbook (many
# "prices" use dots as thousands separators, e.g. "43.011.17", and the
# percentages keep their padding spaces, e.g. "  -4.84%  "). Some of the
# new values happen to look like ordinary decimal numbers (e.g. "6.24"),
# and Excel would normally auto-convert those to numeric cells as soon as
# we assign them. To keep the data as text - matching the workbook's
# existing format - we temporarily force Text number formatting for those
# cells, then restore the default style so we don't leave any formatting
# artifacts behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '42.919.40' },
    @{ Cell = 'E2'; Value = '  -4.71%  ' },
    @{ Cell = 'D3'; Value = '2.217.46' },
    @{ Cell = 'E3'; Value = '  -6.06%  ' },
    @{ Cell = 'D5'; Value = '317.61' },
    @{ Cell = 'E5'; Value = '  +2.26%  ' },
    @{ Cell = 'D6'; Value = '97.52' },
    @{ Cell = 'E6'; Value = '  -9.77%  ' },
    @{ Cell = 'D7'; Value = '0.578' },
    @{ Cell = 'E7'; Value = '  -7.86%  ' },
    @{ Cell = 'E8'; Value = '  +0.02%  ' },
    @{ Cell = 'D9'; Value = '0.561' },
    @{ Cell = 'E9'; Value = '  -8.77%  ' },
    @{ Cell = 'D10'; Value = '36.36' },
    @{ Cell = 'E10'; Value = '  -11.60%  ' },
    @{ Cell = 'D11'; Value = '54.28' },
    @{ Cell = 'E11'; Value = '  -2.15%  ' },
    @{ Cell = 'D12'; Value = '0.0825' },
    @{ Cell = 'E12'; Value = '  -10.13%  ' },
    @{ Cell = 'D13'; Value = '7.70' },
    @{ Cell = 'E13'; Value = '  -8.77%  ' },
    @{ Cell = 'E14'; Value = '  -4.11%  ' },
    @{ Cell = 'D15'; Value = '0.866' },
    @{ Cell = 'E15'; Value = '  -11.87%  ' },
    @{ Cell = 'D16'; Value = '2.556.88' },
    @{ Cell = 'E16'; Value = '  -5.87%  ' },
    @{ Cell = 'D17'; Value = '14.01' },
    @{ Cell = 'E17'; Value = '  -8.31%  ' },
    @{ Cell = 'D18'; Value = '2.202.52' },
    @{ Cell = 'E18'; Value = '  -6.71%  ' },
    @{ Cell = 'D19'; Value = '42.894.97' },
    @{ Cell = 'E19'; Value = '  -4.95%  ' },
    @{ Cell = 'D20'; Value = '14.58' },
    @{ Cell = 'E20'; Value = '  +3.39%  ' },
    @{ Cell = 'D21'; Value = '0.0₃0958' },
    @{ Cell = 'E21'; Value = '  -9.80%  ' },
    @{ Cell = 'D22'; Value = '6.43' },
    @{ Cell = 'E22'; Value = '  -12.17%  ' },
    @{ Cell = 'D23'; Value = '64.98' },
    @{ Cell = 'E23'; Value = '  -11.12%  ' },
    @{ Cell = 'D24'; Value = '3.18' },
    @{ Cell = 'E24'; Value = '  -8.95%  ' },
    @{ Cell = 'D25'; Value = '236.15' },
    @{ Cell = 'E25'; Value = '  -9.13%  ' },
    @{ Cell = 'E26'; Value = '  -8.71%  ' },
    @{ Cell = 'E27'; Value = '  -0.11%  ' },
    @{ Cell = 'D28'; Value = '10.10' },
    @{ Cell = 'E28'; Value = '  -9.56%  ' },
    @{ Cell = 'D30'; Value = '6.24' },
    @{ Cell = 'E30'; Value = '  -14.48%  ' },
    @{ Cell = 'B31'; Value = 'EthereumClassic' },
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' },
    @{ Cell = 'D31'; Value = '20.45' },
    @{ Cell = 'E31'; Value = '  -8.34%  ' },
    @{ Cell = 'B32'; Value = 'Hedera' },
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Cell = 'D32'; Value = '0.0881' },
    @{ Cell = 'E32'; Value = '  -8.77%  ' },
    @{ Cell = 'B33'; Value = 'Monero' },
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Cell = 'D33'; Value = '155.93' },
    @{ Cell = 'E33'; Value = '  -7.44%  ' },
    @{ Cell = 'B34'; Value = 'InjectiveProtocol' },
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' },
    @{ Cell = 'D34'; Value = '33.59' },
    @{ Cell = 'E34'; Value = '  -11.07%  ' },
    @{ Cell = 'E35'; Value = '  -5.62%  ' },
    @{ Cell = 'D36'; Value = '3.29' },
    @{ Cell = 'E36'; Value = '  +10.30%  ' },
    @{ Cell = 'D37'; Value = '2.01' },
    @{ Cell = 'E37'; Value = '  +15.77%  ' },
    @{ Cell = 'E38'; Value = '  -6.39%  ' },
    @{ Cell = 'D39'; Value = '4.44' },
    @{ Cell = 'E39'; Value = '  -7.68%  ' },
    @{ Cell = 'D40'; Value = '0.103' },
    @{ Cell = 'E40'; Value = '  -12.33%  ' },
    @{ Cell = 'D41'; Value = '3.67' },
    @{ Cell = 'E41'; Value = '  -6.74%  ' },
    @{ Cell = 'D42'; Value = '0.0322' },
    @{ Cell = 'E42'; Value = '  -9.02%  ' },
    @{ Cell = 'D43'; Value = '1.879.36' },
    @{ Cell = 'E43'; Value = '  +11.68%  ' },
    @{ Cell = 'E44'; Value = '  +0.24%  ' },
    @{ Cell = 'E45'; Value = '  -5.84%  ' },
    @{ Cell = 'D46'; Value = '88.04' },
    @{ Cell = 'E46'; Value = '  -11.24%  ' },
    @{ Cell = 'D47'; Value = '0.206' },
    @{ Cell = 'E47'; Value = '  -11.19%  ' },
    @{ Cell = 'D48'; Value = '5.47' },
    @{ Cell = 'E48'; Value = '  -0.59%  ' },
    @{ Cell = 'D49'; Value = '77.91' },
    @{ Cell = 'E49'; Value = '  -4.21%  ' },
    @{ Cell = 'D50'; Value = '60.31' },
    @{ Cell = 'E50'; Value = '  -13.30%  ' },
    @{ Cell = 'D51'; Value = '8.67' },
    @{ Cell = 'E51'; Value = '  -5.54%  ' }
)

# Matches plain decimal numbers (optionally signed) with a single decimal
# point, e.g. "6.24", "0.0825", "-3.1" - the only kind of text Excel will
# silently reinterpret as a number on assignment. Values with two dots
# (European-style thousands grouping like "43.011.17"), percent signs,
# padding spaces, or letters are left alone.
$numericLike = '^\s*[+-]?(\d+\.?\d*|\.\d+)\s*$'

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $text = $update.Value

    if ($text -match $numericLike) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}
